$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 336
$ws1.Range("F5").Value = 4954
$ws1.Range("F9").Value = 751
$ws1.Range("F11").Value = 3

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 28

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 336
$ws4.Range("F5").Value = 4954
$ws4.Range("F9").Value = 751
$ws4.Range("F10").Value = 28
$ws4.Range("F12").Value = 3
